$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 29: carries forward the date from B21 (10/10) into column A,
# a new date 12/5 into column B, and the period-79 description into column C.
$ws.Range("A29").Value = "10/10"
$ws.Range("B29").Value = "12/5"
$ws.Range("C29").Value = "第79期 第七代坐騎"

# New row 30: carries forward the date from B22 (10/17) into column A,
# a new date 12/12 into column B, and the period-80 description into column C.
$ws.Range("A30").Value = "10/17"
$ws.Range("B30").Value = "12/12"
$ws.Range("C30").Value = "第80期 祕寶 開放區域 赤潮狂途 祕寶效果: 戰術進階傷害提高1534930 (11051496)"

# Scroll the view down and move the active selection, matching the
# author's on-screen state when they saved (new bottom row selected).
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("C31").Select()
